$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the inventory count for E2 (8 -> 5)
$ws.Range("E2").Value = 5

# Move the active selection on the sheet to E8 (was D3)
$ws.Range("E8").Select()
